# Apply weekly reshuffle of Chirimoya price rows (rows 2-20, columns D:T)
# Each row is re-populated with the data of another row per the source diff,
# effectively reordering the dataset by date without touching static columns A-C, E-K.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44524
$ws.Range("M2").Value = 200
$ws.Range("N2").Value = 23000
$ws.Range("O2").Value = 24000
$ws.Range("P2").Value = 23500
$ws.Range("S2").Value = 1958
$ws.Range("D3").Value = 44860
$ws.Range("R3").Value = "Provincia de Limarí"
$ws.Range("D4").Value = 44167
$ws.Range("N4").Value = 18000
$ws.Range("O4").Value = 19000
$ws.Range("P4").Value = 18500
$ws.Range("S4").Value = 1423
$ws.Range("D5").Value = 44846
$ws.Range("L5").Value = "Primera"
$ws.Range("N5").Value = 24000
$ws.Range("O5").Value = 25000
$ws.Range("P5").Value = 24500
$ws.Range("Q5").Value = "$/caja 12 kilos"
$ws.Range("S5").Value = 2042
$ws.Range("T5").Value = 12
$ws.Range("D6").Value = 44846
$ws.Range("L6").Value = "Segunda"
$ws.Range("M6").Value = 100
$ws.Range("D7").Value = 44874
$ws.Range("L7").Value = "Segunda"
$ws.Range("M7").Value = 250
$ws.Range("N7").Value = 22000
$ws.Range("O7").Value = 23000
$ws.Range("P7").Value = 22500
$ws.Range("Q7").Value = "$/caja 12 kilos"
$ws.Range("S7").Value = 1875
$ws.Range("T7").Value = 12
$ws.Range("D9").Value = 44489
$ws.Range("L9").Value = "Primera"
$ws.Range("N9").Value = 24000
$ws.Range("O9").Value = 25000
$ws.Range("P9").Value = 24500
$ws.Range("S9").Value = 2042
$ws.Range("D10").Value = 44160
$ws.Range("L10").Value = "Segunda"
$ws.Range("M10").Value = 200
$ws.Range("N10").Value = 19000
$ws.Range("O10").Value = 20000
$ws.Range("P10").Value = 19500
$ws.Range("Q10").Value = "$/caja 13 kilos"
$ws.Range("S10").Value = 1500
$ws.Range("T10").Value = 13
$ws.Range("D11").Value = 44811
$ws.Range("L11").Value = "Primera"
$ws.Range("N11").Value = 29000
$ws.Range("O11").Value = 30000
$ws.Range("P11").Value = 29500
$ws.Range("S11").Value = 2458
$ws.Range("D12").Value = 44881
$ws.Range("N12").Value = 22000
$ws.Range("O12").Value = 23000
$ws.Range("P12").Value = 22500
$ws.Range("S12").Value = 1875
$ws.Range("D13").Value = 44475
$ws.Range("L13").Value = "Especial"
$ws.Range("M13").Value = 200
$ws.Range("N13").Value = 32000
$ws.Range("O13").Value = 33000
$ws.Range("P13").Value = 32500
$ws.Range("S13").Value = 2708
$ws.Range("D14").Value = 44482
$ws.Range("L14").Value = "Primera"
$ws.Range("M14").Value = 160
$ws.Range("N14").Value = 25000
$ws.Range("O14").Value = 26000
$ws.Range("P14").Value = 25500
$ws.Range("S14").Value = 2125
$ws.Range("D15").Value = 44783
$ws.Range("L15").Value = "Tercera"
$ws.Range("M15").Value = 100
$ws.Range("N15").Value = 27000
$ws.Range("O15").Value = 28000
$ws.Range("P15").Value = 27500
$ws.Range("S15").Value = 2292
$ws.Range("D16").Value = 44441
$ws.Range("M16").Value = 100
$ws.Range("N16").Value = 29000
$ws.Range("O16").Value = 30000
$ws.Range("P16").Value = 29500
$ws.Range("R16").Value = "Región de Coquimbo"
$ws.Range("S16").Value = 2458
$ws.Range("D17").Value = 44776
$ws.Range("M17").Value = 160
$ws.Range("N17").Value = 29000
$ws.Range("O17").Value = 30000
$ws.Range("P17").Value = 29500
$ws.Range("Q17").Value = "$/caja 10 kilos"
$ws.Range("S17").Value = 2950
$ws.Range("T17").Value = 10
$ws.Range("D18").Value = 44839
$ws.Range("L18").Value = "Segunda"
$ws.Range("M18").Value = 160
$ws.Range("N18").Value = 26000
$ws.Range("O18").Value = 27000
$ws.Range("P18").Value = 26500
$ws.Range("Q18").Value = "$/caja 12 kilos"
$ws.Range("S18").Value = 2208
$ws.Range("D19").Value = 44468
$ws.Range("L19").Value = "Primera"
$ws.Range("M19").Value = 200
$ws.Range("N19").Value = 29000
$ws.Range("O19").Value = 30000
$ws.Range("P19").Value = 29500
$ws.Range("Q19").Value = "$/bandeja 10 kilos"
$ws.Range("S19").Value = 2950
$ws.Range("T19").Value = 10
$ws.Range("D20").Value = 44545
$ws.Range("L20").Value = "Primera"
$ws.Range("M20").Value = 200
$ws.Range("N20").Value = 23000
$ws.Range("O20").Value = 24000
$ws.Range("P20").Value = 23500
$ws.Range("Q20").Value = "$/bandeja 12 kilos"
$ws.Range("S20").Value = 1958
